# Update gh-pages to output generated at 456a3b4
#
# Sheet "展览" (1): bump "想去人数" (F) counters on a few rows, and append a
# new row (11) for the "肥西·星域动漫游戏嘉年华" event.
# Sheet "演出" (2): bump F2 counter by 1.
# Sheet "全部类型" (4): bump the same F counters (this sheet aggregates all
# events), and insert the new "肥西·星域动漫游戏嘉年华" event as row 12
# (between the existing rows), pushing the old row 12 down to row 13.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Setting a plain string that looks like a date (e.g. "2024-09-16")
    # gets auto-converted to a date serial by Excel's normal typed-value
    # coercion. Force the cell to text first, then strip the resulting
    # NumberFormat override again so the cell keeps the workbook's default
    # (unstyled) look, matching the other text cells in the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Copy-CellFormat {
    param($srcCell, $dstCell)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(3, 6).Value = 101
$ws1.Cells.Item(4, 6).Value = 496
$ws1.Cells.Item(5, 6).Value = 4916
$ws1.Cells.Item(9, 6).Value = 747
$ws1.Cells.Item(10, 6).Value = 232

Copy-CellFormat $ws1.Cells.Item(10, 1) $ws1.Cells.Item(11, 1)

$ws1.Cells.Item(11, 1).Value = 10
Set-TextValue $ws1.Cells.Item(11, 2) "2024-09-16"
$ws1.Cells.Item(11, 3).Value = "肥西·星域动漫游戏嘉年华"
$ws1.Cells.Item(11, 4).Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws1.Cells.Item(11, 5).Value = "2024.09.16 10:00-09.16 17:00"
$ws1.Cells.Item(11, 6).Value = 1
$ws1.Cells.Item(11, 7).Value = 45
$ws1.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws1.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(2, 6).Value = 26

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(3, 6).Value = 101
$ws4.Cells.Item(4, 6).Value = 496
$ws4.Cells.Item(5, 6).Value = 4916
$ws4.Cells.Item(9, 6).Value = 747
$ws4.Cells.Item(10, 6).Value = 26
$ws4.Cells.Item(11, 6).Value = 232

# Insert a new row 12 (shifts the old row 12 "四月是你的谎言" down to 13),
# and renumber the leading index column to keep it sequential.
$ws4.Rows.Item(12).Insert()

Copy-CellFormat $ws4.Cells.Item(11, 1) $ws4.Cells.Item(12, 1)

$ws4.Cells.Item(12, 1).Value = 11
$ws4.Cells.Item(13, 1).Value = 12

Set-TextValue $ws4.Cells.Item(12, 2) "2024-09-16"
$ws4.Cells.Item(12, 3).Value = "肥西·星域动漫游戏嘉年华"
$ws4.Cells.Item(12, 4).Value = "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)"
$ws4.Cells.Item(12, 5).Value = "2024.09.16 10:00-09.16 17:00"
$ws4.Cells.Item(12, 6).Value = 1
$ws4.Cells.Item(12, 7).Value = 45
$ws4.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90489"
$ws4.Cells.Item(12, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

Write-Host "edit applied"
